$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2021-10-05"

# Row 3 (January) - 2021 column updates
$ws.Range("U3").Value = 202
$ws.Range("V3").Value = 0.0691

# Row 12 (October) - update label and figures for new day's data
$ws.Range("A12").Value = "October (through 10-05)"
$ws.Range("B12").Value = 1
$ws.Range("D12").Value = 0.2
$ws.Range("D12").NumberFormat = "0.0%"
$ws.Range("F12").Value = 7
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 0.1818
$ws.Range("L12").Value = 12
$ws.Range("O12").Value = 4
$ws.Range("R12").Value = 29
$ws.Range("U12").Value = 35

# Row 13 (Total) - updated totals
$ws.Range("B13").Value = 31
$ws.Range("D13").Value = 0.1342
$ws.Range("F13").Value = 390
$ws.Range("G13").Value = 0.1055
$ws.Range("I13").Value = 586
$ws.Range("J13").Value = 0.0815
$ws.Range("L13").Value = 499
$ws.Range("M13").Value = 0.1089
$ws.Range("O13").Value = 383
$ws.Range("P13").Value = 0.1009
$ws.Range("R13").Value = 877
$ws.Range("S13").Value = 0.057
$ws.Range("U13").Value = 1205
$ws.Range("V13").Value = 0.0615
